# Working IEEE models and ORPD.
# Add an "I_lim_A" (line current limit) column to the "lines" sheet, right
# after v_nom_kv and before length_km, and drop the stray number-format
# style that had been applied to gens!C3.

$wb = $excel.ActiveWorkbook

# --- 1. "lines" sheet: insert the new I_lim_A column at column C ---------
$ws = $wb.Worksheets.Item("lines")

# Shift the existing columns C:I (length_km .. is_pu) one place to the
# right, into D:J, without touching the <cols> width definitions (work
# right-to-left so we never clobber a column before it has been copied).
for ($r = 1; $r -le 4; $r++) {
    for ($c = 9; $c -ge 3; $c--) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r, $c + 1)
        $dst.Value = $src.Value2
    }
}

# Populate the freed-up column C with the new field.
$ws.Cells.Item(1, 3).Value = "I_lim_A"
$ws.Cells.Item(2, 3).Value = 1000
$ws.Cells.Item(3, 3).Value = 1000
$ws.Cells.Item(4, 3).Value = 1000

# --- 2. "gens" sheet: clear the leftover applyNumberFormat style on C3 ---
$gens = $wb.Worksheets.Item("gens")
$gens.Cells.Item(3, 3).ClearFormats()

# --- 3. Leave the UI focused on the "lines" sheet / new column ----------
[void]$ws.Range("C4").Select()
$ws.Activate()
